$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.931.64'
$ws.Range("E2").Value = '  +0.32%  '

$ws.Range("D3").Value = '2.282.47'
$ws.Range("E3").Value = '  +0.29%  '

$ws.Range("D5").Value = '''249.67'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.66%  '

$ws.Range("D6").Value = '''0.634'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.14%  '

$ws.Range("D7").Value = '''78.96'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +10.26%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").Value = '''0.653'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.33%  '

$ws.Range("D10").Value = '''41.02'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.67%  '

$ws.Range("D11").Value = '''0.0974'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.71%  '

$ws.Range("D12").Value = '''7.33'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.26%  '

$ws.Range("E13").Value = '  +0.11%  '

$ws.Range("D14").Value = '2.621.31'
$ws.Range("E14").Value = '  +0.12%  '

$ws.Range("D15").Value = '''15.04'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.63%  '

$ws.Range("D16").Value = '''0.868'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.73%  '

$ws.Range("D17").Value = '2.287.53'
$ws.Range("E17").Value = '  +0.84%  '

$ws.Range("D18").Value = '42.825.91'
$ws.Range("E18").Value = '  +0.14%  '

$ws.Range("D19").Value = '0.0₃0992'
$ws.Range("E19").Value = '  -2.47%  '

$ws.Range("D20").Value = '''6.22'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.57%  '

$ws.Range("D21").Value = '''72.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.72%  '

$ws.Range("D22").Value = '''234.10'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.15%  '

$ws.Range("D23").Value = '''2.17'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.86%  '

$ws.Range("E24").Value = '  -2.81%  '

$ws.Range("E25").Value = '  -0.14%  '

$ws.Range("D26").Value = '''11.36'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.23%  '

$ws.Range("E27").Value = '  -4.07%  '

$ws.Range("D28").Value = '''2.17'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.31%  '

$ws.Range("D29").Value = '''167.61'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.46%  '

$ws.Range("E30").Value = '  -2.08%  '

$ws.Range("D31").Value = '''6.47'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.07%  '

$ws.Range("E32").Value = '  +7.04%  '

$ws.Range("E33").Value = '  -5.29%  '

$ws.Range("D34").Value = '''30.20'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.58%  '

$ws.Range("D35").Value = '''0.127'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.65%  '

$ws.Range("D36").Value = '''4.57'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.66%  '

$ws.Range("D37").Value = '''4.78'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.52%  '

$ws.Range("D38").Value = '''0.0305'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.90%  '

$ws.Range("D39").Value = '''13.78'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.32%  '

$ws.Range("E40").Value = '  -2.38%  '

$ws.Range("E41").Value = '  +0.14%  '

$ws.Range("D42").Value = '''112.95'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +17.47%  '

$ws.Range("D43").Value = '''0.206'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.41%  '

$ws.Range("D44").Value = '''61.23'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.75%  '

$ws.Range("E45").Value = '  -2.85%  '

$ws.Range("E46").Value = '  -0.70%  '

$ws.Range("D49").Value = '''1.15'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.75%  '

$ws.Range("D50").Value = '''1.16'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.36%  '

$ws.Range("D51").Value = '''4.24'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.03%  '

$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").Value = '''4.62'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.47%  '

$ws.Range("B48").Value = 'BinanceUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D48").Value = '''1.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.05%  '
